# Scheduled market-price refresh: update currentAveragePrice* / Leve*Profit* columns
# pulled from the Universalis API snapshot for each Leve row across the DoH sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18: You Grow, Girl (Leve Item ID 5471)
$ws.Range("H18").Value = 997.5
$ws.Range("J18").Value = 990
$ws.Range("L18").Value = 990
$ws.Range("N18").Value = -1558

# Row 40: Stuck in the Moment (Leve Item ID 5505)
$ws.Range("H40").Value = 4249.8335
$ws.Range("I40").Value = 2944.3333
$ws.Range("K40").Value = 2944.3333
$ws.Range("M40").Value = -2769.3333

# Row 43: Growing Is Knowing (Leve Item ID 5472)
$ws.Range("H43").Value = 2653.625
$ws.Range("J43").Value = 2045.6
$ws.Range("L43").Value = 2045.6
$ws.Range("N43").Value = -2183.6

# Row 74: Adhesive of Antipathy (Leve Item ID 5507)
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").ClearContents()

# Row 76: Warding Off Temptation (Leve Item ID 12602)
$ws.Range("H76").Value = 3585.2856
$ws.Range("I76").Value = 3499.4
$ws.Range("K76").Value = 3499.4
$ws.Range("M76").Value = -3184.4

# Row 77: It's Gonna Grow Back (L) (Leve Item ID 5507)
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").ClearContents()

# Row 79: The Garden of Arcane Delights (L) (Leve Item ID 12602)
$ws.Range("H79").Value = 3585.2856
$ws.Range("I79").Value = 3499.4
$ws.Range("K79").Value = 3499.4
$ws.Range("M79").Value = -2407.4

# Row 125: Body over Mind (Leve Item ID 36228)
$ws.Range("H125").Value = 2814.8462
$ws.Range("I125").Value = 2374.5715
$ws.Range("J125").Value = 3328.5
$ws.Range("K125").Value = 21371.1435
$ws.Range("L125").Value = 29956.5
$ws.Range("M125").Value = -18911.1435
$ws.Range("N125").Value = -34876.5

# Row 131: Mindful Study (Leve Item ID 36108)
$ws.Range("H131").Value = 369.5
$ws.Range("I131").Value = 369.5
$ws.Range("K131").Value = 1108.5
$ws.Range("M131").Value = 3931.5

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust (Leve Item ID 44147)
$ws.Range("H32").Value = 3696
$ws.Range("I32").Value = 2748.25
$ws.Range("J32").Value = 7487
$ws.Range("K32").Value = 2748.25
$ws.Range("L32").Value = 7487
$ws.Range("M32").Value = -2461.25
$ws.Range("N32").Value = -8061

# Row 61: Dealing with the Tough Stuff (Leve Item ID 43999)
$ws.Range("H61").Value = 27507
$ws.Range("J61").Value = 27507
$ws.Range("L61").Value = 27507
$ws.Range("N61").Value = -27931

# Row 74: As the Bolt Flies (Leve Item ID 44000)
$ws.Range("H74").Value = 3497.6667
$ws.Range("I74").Value = 3497.6667
$ws.Range("K74").Value = 3497.6667
$ws.Range("M74").Value = -2623.6667

# Row 77: Heavy Metal Banned (L) (Leve Item ID 44000)
$ws.Range("H77").Value = 3497.6667
$ws.Range("I77").Value = 3497.6667
$ws.Range("K77").Value = 17488.3335
$ws.Range("M77").Value = -13120.3335

# Row 97: Ore for Me (Leve Item ID 19941)
$ws.Range("H97").Value = 1313.28
$ws.Range("I97").Value = 1074.6666
$ws.Range("J97").Value = 1926.8572
$ws.Range("K97").Value = 1074.6666
$ws.Range("L97").Value = 1926.8572
$ws.Range("M97").Value = -578.6666
$ws.Range("N97").Value = -2918.8572

# Row 132: Don't Bore Me, Ore Me (Leve Item ID 43997)
$ws.Range("H132").Value = 1881.0667
$ws.Range("I132").Value = 1729.7142
$ws.Range("K132").Value = 5189.142599999999
$ws.Range("M132").Value = -2659.142599999999

# Row 136: Metal with Mettle (Leve Item ID 43999)
$ws.Range("H136").Value = 27507
$ws.Range("J136").Value = 27507
$ws.Range("L136").Value = 82521
$ws.Range("N136").Value = -87621

$ws = $wb.Worksheets.Item("CRP")
# Row 28: Militia on My Mind (Leve Item ID 18348)
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()

# Row 31: Wall Not Found (Leve Item ID 44023)
$ws.Range("H31").Value = 1749
$ws.Range("I31").Value = 1749
$ws.Range("K31").Value = 1749
$ws.Range("M31").Value = -1454

# Row 34: Armoires of the Rich and Famous (Leve Item ID 44023)
$ws.Range("H34").Value = 1749
$ws.Range("I34").Value = 1749
$ws.Range("K34").Value = 1749
$ws.Range("M34").Value = -1547

# Row 58: You Do the Heavy Lifting (Leve Item ID 44021)
$ws.Range("H58").Value = 3420.6
$ws.Range("I58").Value = 2887
$ws.Range("K58").Value = 2887
$ws.Range("M58").Value = -2684

# Row 80: The Long Armillae of the Law (Leve Item ID 12015)
$ws.Range("H80").Value = 15000
$ws.Range("I80").Value = 10000
$ws.Range("K80").Value = 10000
$ws.Range("M80").Value = -8877

# Row 83: Wooden Ambitions (L) (Leve Item ID 12015)
$ws.Range("H83").Value = 15000
$ws.Range("I83").Value = 10000
$ws.Range("K83").Value = 30000
$ws.Range("M83").Value = -24384

# Row 108: Just Starting Out (Leve Item ID 27087)
$ws.Range("H108").Value = 30000
$ws.Range("J108").Value = 30000
$ws.Range("L108").Value = 30000
$ws.Range("N108").Value = -37680

# Row 132: Hull Lotta Damage (Leve Item ID 44019)
$ws.Range("H132").Value = 2446.3635
$ws.Range("I132").Value = 2446.3635
$ws.Range("K132").Value = 7339.0905
$ws.Range("M132").Value = -4809.0905

# Row 136: Turali Quality (Leve Item ID 44021)
$ws.Range("H136").Value = 3420.6
$ws.Range("I136").Value = 2887
$ws.Range("K136").Value = 8661
$ws.Range("M136").Value = -6111

$ws = $wb.Worksheets.Item("GSM")
# Row 15: The Tusk at Hand (Leve Item ID 12018)
$ws.Range("H15").Value = 25000
$ws.Range("J15").Value = 25000
$ws.Range("L15").Value = 25000
$ws.Range("N15").Value = -25576

# Row 81: The Grander Temple (Leve Item ID 12018)
$ws.Range("H81").Value = 25000
$ws.Range("J81").Value = 25000
$ws.Range("L81").Value = 25000
$ws.Range("N81").Value = -26996

# Row 84: Man with a Dragon Earring (L) (Leve Item ID 12018)
$ws.Range("H84").Value = 25000
$ws.Range("J84").Value = 25000
$ws.Range("L84").Value = 75000
$ws.Range("N84").Value = -84984

# Row 132: On Board for Lar (Leve Item ID 44008)
$ws.Range("H132").Value = 2082.2
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 2082.2
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 6246.599999999999
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -11306.6

$ws = $wb.Worksheets.Item("LTW")
# Row 46: Supply Side Logic (Leve Item ID 5282)
$ws.Range("H46").Value = 4132.243
$ws.Range("I46").Value = 3827.8
$ws.Range("J46").Value = 4766.5
$ws.Range("K46").Value = 3827.8
$ws.Range("L46").Value = 4766.5
$ws.Range("M46").Value = -3639.8
$ws.Range("N46").Value = -5142.5

# Row 61: Spelling Me Softly (Leve Item ID 27740)
$ws.Range("H61").Value = 3199
$ws.Range("I61").Value = 2978.8
$ws.Range("J61").Value = 3749.5
$ws.Range("K61").Value = 2978.8
$ws.Range("L61").Value = 3749.5
$ws.Range("M61").Value = -2776.8
$ws.Range("N61").Value = -4153.5

# Row 68: You Could Say It's a Moving Target (Leve Item ID 12563)
$ws.Range("H68").Value = 1424.3334
$ws.Range("I68").Value = 1541.8572
$ws.Range("J68").Value = 1259.8
$ws.Range("K68").Value = 1541.8572
$ws.Range("L68").Value = 1259.8
$ws.Range("M68").Value = -792.8571999999999
$ws.Range("N68").Value = -2757.8

# Row 71: They Call It Bloody Mary (L) (Leve Item ID 12563)
$ws.Range("H71").Value = 1424.3334
$ws.Range("I71").Value = 1541.8572
$ws.Range("J71").Value = 1259.8
$ws.Range("K71").Value = 7709.286
$ws.Range("L71").Value = 6299
$ws.Range("M71").Value = -3965.286
$ws.Range("N71").Value = -13787

# Row 113: Peace in Rest (Leve Item ID 27740)
$ws.Range("H113").Value = 3199
$ws.Range("I113").Value = 2978.8
$ws.Range("J113").Value = 3749.5
$ws.Range("K113").Value = 2978.8
$ws.Range("L113").Value = 3749.5
$ws.Range("M113").Value = -808.8000000000002
$ws.Range("N113").Value = -8089.5

# Row 131: For What Was Gleaned (Leve Item ID 35466)
$ws.Range("H131").Value = 55999.8
$ws.Range("J131").Value = 55999.8
$ws.Range("L131").Value = 55999.8
$ws.Range("N131").Value = -66079.8

# Row 132: Tenets of Tanning (Leve Item ID 44058)
$ws.Range("H132").Value = 3973
$ws.Range("I132").Value = 3858.889
$ws.Range("K132").Value = 11576.667
$ws.Range("M132").Value = -9046.667000000001
